# Applies the "Add files via upload" edit: appends a new sub-topic
# (new prospective technologies exercise) under item 6, and makes
# cell C39 bold to match the other sub-heading cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix formatting on C39: make it bold + wrap, matching the other
#     "sub-heading" style used for C19 / C27 / C43, and give it the
#     taller row height used by wrapped rows.
$ws.Range("C39").Font.Bold = $true
$ws.Range("C39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 34.9

# --- Append the new rows (44-49) describing the "new prospective
#     technologies" sub-exercise.
$ws.Range("C44").Value = 'From "Files for Later\SubRES_TMPL" copy SubRES_NEW_PP.xlsx file into \SubRES_TMPL folder'

$ws.Range("C45").Value = "Upload the new file (remember to put it into \SubRES_TMPL folder folder)"
$ws.Range("C45").WrapText = $true

$ws.Range("C46").Value = "What has changed in the Navigator module?"
$ws.Range("C46").WrapText = $true

$ws.Range("C47").Value = "Look at the attributes in the top FI_T table. Which attributes are new?"

$ws.Range("C48").Value = "Synchronize and run the model up to the year 2050. Check up the results."

$ws.Range("C49").Value = "Add a new process of your choice. Define the process and commodities if necessary. Set values for attributes. Upload, synchronize and run the model."
$ws.Range("C49").WrapText = $true
$ws.Rows.Item(49).RowHeight = 34.9

# --- Update the view so the newly added rows are visible, matching
#     the author's final selection/scroll position.
$ws.Range("C39").Select()
$excel.ActiveWindow.ScrollRow = 28
